# [FEATURE] Compensacion Coelsa + cambios regresion
#
# Inserts a new "Carga Saldos" worksheet between "Users" and "DatosNCD",
# populates it with the Cuenta/Importe/Concepto header row plus one data
# row, sizes column A, and makes the new sheet the active/selected tab
# (moving tabSelected off of "Users").

$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")

# Adding "before" the sheet that currently follows Users (DatosNCD) would
# also work, but adding "after" Users directly is more robust to ordering.
$newSheet = $wb.Worksheets.Add($null, $usersSheet)
$newSheet.Name = "Carga Saldos"

# Header row
$newSheet.Range("A1").Value = "Cuenta"
$newSheet.Range("B1").Value = "Importe "
$newSheet.Range("C1").Value = "Concepto"

# Data row
$newSheet.Range("A2").Value = 10010656046
$newSheet.Range("B2").Value = 100
$newSheet.Range("C2").Value = "18602AME"

# Column A width -> stored sheet width of 14 (Excel adds ~0.8333 padding
# to the character-based ColumnWidth when serialising to the <col> width).
$newSheet.Columns.Item(1).ColumnWidth = 13.166666666666666

# Make the new sheet the active tab / selected cell, matching the
# activeTab change (2 -> 3) and the move of tabSelected from Users.
$newSheet.Activate()
$newSheet.Range("C2").Select()
